# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F9").Value = 29
    $ws.Range("F15").Value = 940
    $ws.Range("F18").Value = 421
    $ws.Range("F23").Value = 1281
    $ws.Range("F24").Value = 2934
    $ws.Range("F28").Value = 67
    $ws.Range("F33").Value = 272
    $ws.Range("F36").Value = 615
}
